$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@(16, '1047399189', 'JOSE DEL CARMEN CUADRO SAENZ', '2206', 22916, 737717)
    ,@(17, '1047399189', 'JOSE DEL CARMEN CUADRO SAENZ', '2205', 31249, 737717)
    ,@(18, '1047399189', 'JOSE DEL CARMEN CUADRO SAENZ', '2204', 31249, 737717)
    ,@(19, '1047399189', 'JOSE DEL CARMEN CUADRO SAENZ', '2203', 31249, 737717)
    ,@(20, '1047399189', 'JOSE DEL CARMEN CUADRO SAENZ', '2202', 31249, 737717)
    ,@(21, '1047399189', 'JOSE DEL CARMEN CUADRO SAENZ', '2201', 31249, 737717)
    ,@(22, '1047399189', 'JOSE DEL CARMEN CUADRO SAENZ', '2112', 31249, 737717)
    ,@(23, '1047399189', 'JOSE DEL CARMEN CUADRO SAENZ', '2111', 31249, 737717)
    ,@(24, '1047399189', 'JOSE DEL CARMEN CUADRO SAENZ', '2110', 31249, 737717)
    ,@(25, '1047399189', 'JOSE DEL CARMEN CUADRO SAENZ', '2109', 31249, 737717)
    ,@(26, '1047399189', 'JOSE DEL CARMEN CUADRO SAENZ', '2108', 31249, 737717)
    ,@(27, '1047399189', 'JOSE DEL CARMEN CUADRO SAENZ', '2107', 31249, 737717)
    ,@(28, '1047399189', 'JOSE DEL CARMEN CUADRO SAENZ', '2106', 31249, 737717)
    ,@(29, '1047399189', 'JOSE DEL CARMEN CUADRO SAENZ', '2105', 31249, 737717)
    ,@(30, '1047399189', 'JOSE DEL CARMEN CUADRO SAENZ', '2104', 31249, 737717)
    ,@(31, '1047399189', 'JOSE DEL CARMEN CUADRO SAENZ', '2103', 31249, 737717)
    ,@(32, '1047399189', 'JOSE DEL CARMEN CUADRO SAENZ', '2102', 31249, 737717)
    ,@(33, '1047399189', 'JOSE DEL CARMEN CUADRO SAENZ', '2101', 31249, 737717)
    ,@(34, '1047399189', 'JOSE DEL CARMEN CUADRO SAENZ', '2012', 31249, 737717)
    ,@(35, '1047399189', 'JOSE DEL CARMEN CUADRO SAENZ', '2011', 31249, 737717)
    ,@(36, '1047399189', 'JOSE DEL CARMEN CUADRO SAENZ', '2010', 31249, 737717)
    ,@(37, '1047399189', 'JOSE DEL CARMEN CUADRO SAENZ', '2009', 31249, 737717)
    ,@(38, '1047399189', 'JOSE DEL CARMEN CUADRO SAENZ', '2008', 31249, 737717)
    ,@(39, '1047399189', 'JOSE DEL CARMEN CUADRO SAENZ', '2007', 31249, 737717)
    ,@(40, '1047399189', 'JOSE DEL CARMEN CUADRO SAENZ', '2006', 31249, 737717)
    ,@(41, '1047399189', 'JOSE DEL CARMEN CUADRO SAENZ', '2005', 31249, 737717)
    ,@(42, '1047399189', 'JOSE DEL CARMEN CUADRO SAENZ', '2004', 31249, 737717)
    ,@(43, '1047399189', 'JOSE DEL CARMEN CUADRO SAENZ', '2003', 31249, 737717)
    ,@(44, '1047399189', 'JOSE DEL CARMEN CUADRO SAENZ', '2002', 31249, 737717)
    ,@(45, '1047399189', 'JOSE DEL CARMEN CUADRO SAENZ', '2001', 31249, 737717)
    ,@(46, '1047399189', 'JOSE DEL CARMEN CUADRO SAENZ', '1912', 31249, 737717)
    ,@(47, '1047399189', 'JOSE DEL CARMEN CUADRO SAENZ', '1911', 31249, 737717)
    ,@(48, '1047399189', 'JOSE DEL CARMEN CUADRO SAENZ', '1910', 31249, 737717)
    ,@(49, '1047399189', 'JOSE DEL CARMEN CUADRO SAENZ', '1909', 31249, 737717)
    ,@(50, '1047399189', 'JOSE DEL CARMEN CUADRO SAENZ', '1908', 31249, 737717)
    ,@(51, '1047399189', 'JOSE DEL CARMEN CUADRO SAENZ', '1907', 31249, 737717)
    ,@(52, '1047399189', 'JOSE DEL CARMEN CUADRO SAENZ', '1906', 31249, 737717)
    ,@(53, '1047399189', 'JOSE DEL CARMEN CUADRO SAENZ', '1905', 31249, 737717)
    ,@(54, '1047399189', 'JOSE DEL CARMEN CUADRO SAENZ', '1904', 31249, 737717)
    ,@(55, '1047399189', 'JOSE DEL CARMEN CUADRO SAENZ', '1903', 31249, 737717)
    ,@(56, '1047399189', 'JOSE DEL CARMEN CUADRO SAENZ', '1902', 31249, 737717)
    ,@(57, '1047399189', 'JOSE DEL CARMEN CUADRO SAENZ', '1901', 31249, 737717)
    ,@(58, '1047399189', 'JOSE DEL CARMEN CUADRO SAENZ', '1812', 31249, 737717)
    ,@(59, '1047399189', 'JOSE DEL CARMEN CUADRO SAENZ', '1811', 31249, 737717)
    ,@(60, '1047399189', 'JOSE DEL CARMEN CUADRO SAENZ', '1810', 31249, 737717)
    ,@(61, '1047399189', 'JOSE DEL CARMEN CUADRO SAENZ', '1809', 31249, 737717)
    ,@(62, '1047399189', 'JOSE DEL CARMEN CUADRO SAENZ', '1808', 29509, 737717)
    ,@(63, '1047399189', 'JOSE DEL CARMEN CUADRO SAENZ', '1807', 29509, 737717)
    ,@(64, '1047399189', 'JOSE DEL CARMEN CUADRO SAENZ', '1806', 29509, 737717)
    ,@(65, '1047399189', 'JOSE DEL CARMEN CUADRO SAENZ', '1805', 29509, 737717)
    ,@(66, '1047399189', 'JOSE DEL CARMEN CUADRO SAENZ', '1804', 29509, 737717)
    ,@(67, '1047399189', 'JOSE DEL CARMEN CUADRO SAENZ', '1803', 29509, 737717)
    ,@(68, '1047399189', 'JOSE DEL CARMEN CUADRO SAENZ', '1802', 29509, 737717)
    ,@(69, '1047399189', 'JOSE DEL CARMEN CUADRO SAENZ', '1801', 29509, 737717)
    ,@(70, '1047399189', 'JOSE DEL CARMEN CUADRO SAENZ', '1712', 29509, 737717)
    ,@(71, '1047399189', 'JOSE DEL CARMEN CUADRO SAENZ', '1711', 29509, 737717)
    ,@(72, '1047399189', 'JOSE DEL CARMEN CUADRO SAENZ', '1710', 29509, 737717)
    ,@(73, '1047399189', 'JOSE DEL CARMEN CUADRO SAENZ', '1709', 29509, 737717)
    ,@(74, '1047399189', 'JOSE DEL CARMEN CUADRO SAENZ', '1708', 29509, 737717)
    ,@(75, '1143345378', 'RODRIGO JOSE MARIMON UTRIA', '1707', 29509, 737717)
    ,@(76, '1050964013', 'ELKIN FABIAN PAJARO ARELLANO', '1708', 29509, 781242)
    ,@(77, '1050964013', 'ELKIN FABIAN PAJARO ARELLANO', '1707', 7869, 781242)
    ,@(78, '73213446', 'LEONARDO RIVERA SAENZ', '1902', 31249, 737717)
    ,@(79, '73213446', 'LEONARDO RIVERA SAENZ', '1901', 31249, 737717)
    ,@(80, '73213446', 'LEONARDO RIVERA SAENZ', '1812', 31249, 737717)
    ,@(81, '73213446', 'LEONARDO RIVERA SAENZ', '1811', 31249, 737717)
    ,@(82, '73213446', 'LEONARDO RIVERA SAENZ', '1810', 31249, 737717)
    ,@(83, '73213446', 'LEONARDO RIVERA SAENZ', '1809', 31249, 737717)
    ,@(84, '73213446', 'LEONARDO RIVERA SAENZ', '1808', 29509, 737717)
    ,@(85, '73213446', 'LEONARDO RIVERA SAENZ', '1807', 29509, 737717)
    ,@(86, '73213446', 'LEONARDO RIVERA SAENZ', '1806', 29509, 737717)
    ,@(87, '73213446', 'LEONARDO RIVERA SAENZ', '1805', 29509, 737717)
    ,@(88, '73213446', 'LEONARDO RIVERA SAENZ', '1804', 29509, 737717)
    ,@(89, '73213446', 'LEONARDO RIVERA SAENZ', '1803', 29509, 737717)
    ,@(90, '73213446', 'LEONARDO RIVERA SAENZ', '1802', 29509, 737717)
    ,@(91, '73213446', 'LEONARDO RIVERA SAENZ', '1801', 29509, 737717)
    ,@(92, '73213446', 'LEONARDO RIVERA SAENZ', '1712', 29509, 737717)
    ,@(93, '73213446', 'LEONARDO RIVERA SAENZ', '1711', 29509, 737717)
    ,@(94, '73213446', 'LEONARDO RIVERA SAENZ', '1710', 29509, 737717)
    ,@(95, '73213446', 'LEONARDO RIVERA SAENZ', '1709', 29509, 737717)
    ,@(96, '73213446', 'LEONARDO RIVERA SAENZ', '1708', 29509, 737717)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
    $ws.Cells.Item($r, 5).Value = $row[3]
    $ws.Cells.Item($r, 6).Value = $row[4]
    $ws.Cells.Item($r, 7).Value = $row[5]
}
